$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales_Tests")

$ws.Range("A5").Value = "TC_SALE_03"
$ws.Range("B5").Value = "Create Multi-Item Sale"
$ws.Range("A13").Value = "TC_SALE_04"
$ws.Range("B13").Value = "Verify Invoice & Print"
$ws.Range("C13").Value = '1.Verify URL contains "invoice"'

$ws.Range("C5").Value = "1.Click ""Create New Sale"" at ""//a[contains(@class, 'btn-create')]"""
$ws.Range("C6").Value = "2.Type ""faizal"" at ""//select[@name='customerId']"""
$ws.Range("C7").Value = "3.Type ""Dettol (₹80.0)"" at ""(//select[@name='productIds'])[1]"""
$ws.Range("C8").Value = "4.Type ""3"" at ""(//input[@name='quantities'])[1]"""
$ws.Range("C9").Value = "5.Click ""Add Item"" at ""//button[contains(text(), '+ Add Another Item')]"""
$ws.Range("C10").Value = "6.Type ""Pears Soap (₹40.0)"" at ""(//select[@name='productIds'])[2]"""
$ws.Range("C11").Value = "7.Type ""2"" at ""(//input[@name='quantities'])[2]"""
$ws.Range("C12").Value = "8.Click ""Generate Bill"" at ""//button[@type='submit']"""

$ws.Range("C13").Select()
